# Add the new resale-numbers row (2023-06-20 10:33) to the CityResaleNum sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 61

# Columns A (Date) and D (Week) look like a date / a plain number to Excel's
# auto-detection, so force them to be stored as literal text (matching the
# existing rows, which keep these as text values) before assigning.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-20"
$ws.Cells.Item($row, 2).Value = "10:33:11"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 4).Value = "25"

$ws.Cells.Item($row, 5).Value = 122037
$ws.Cells.Item($row, 6).Value = 133735
$ws.Cells.Item($row, 7).Value = 162269
$ws.Cells.Item($row, 8).Value = 133259
$ws.Cells.Item($row, 9).Value = 177288
$ws.Cells.Item($row, 10).Value = 114427
$ws.Cells.Item($row, 11).Value = 201421
$ws.Cells.Item($row, 12).Value = 225175
$ws.Cells.Item($row, 13).Value = 175614
$ws.Cells.Item($row, 14).Value = 103850
$ws.Cells.Item($row, 15).Value = 39187
$ws.Cells.Item($row, 16).Value = 33917
$ws.Cells.Item($row, 17).Value = 51824
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35922
$ws.Cells.Item($row, 20).Value = -1
